$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing values (5) in rows 8, 12 and 27
$ws.Range("C8:E8").Value = 5
$ws.Range("C12:F12").Value = 5
$ws.Range("C27:F27").Value = 5

# Update the frozen pane top-left cell and the active selection
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Application.ActiveWindow.ScrollColumn = 3

$ws.Range("F27").Select()
